$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-10-02 16:46:13", "hatespeech", "def", 20, 0.5007948904606515),
    @("2023-10-02 16:46:13", "hatespeech", "def", 10, 0.6570267269092438),
    @("2023-10-02 16:46:13", "hatespeech", "def", 30, 0.4054775467285684),
    @("2023-10-02 16:46:13", "hatespeech", "def", 40, 0.3301035082478216)
)

$startRow = 34
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}
